# Update "想去人数" (F column) counts on the "展览" and "全部类型" sheets.
# Both sheets carry the same event list, so the same row/value updates
# apply to each.

$wb = $excel.ActiveWorkbook

$updates = @{
    4  = 732
    9  = 447
    23 = 1383
    25 = 316
    38 = 570
    40 = 3490
    41 = 413
    45 = 60
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
